# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to match the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5547
$ws1.Range("G2").Value = 70
$ws1.Range("F3").Value = 621
$ws1.Range("F4").Value = 12510
$ws1.Range("F5").Value = 306
$ws1.Range("F6").Value = 619
$ws1.Range("F7").Value = 191
$ws1.Range("F8").Value = 365
$ws1.Range("F9").Value = 1160
$ws1.Range("F10").Value = 110

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5547
$ws4.Range("G2").Value = 70
$ws4.Range("F3").Value = 621
$ws4.Range("F5").Value = 12510
$ws4.Range("F6").Value = 306
$ws4.Range("F7").Value = 619
$ws4.Range("F8").Value = 191
$ws4.Range("F11").Value = 365
$ws4.Range("F12").Value = 1160
$ws4.Range("F14").Value = 110
